# Apply updates to the event schedule workbook:
# 1) Shift all Time Block dates in column A (rows 2-97) forward by 1 day.
# 2) Update the payload (column B) so that rows for 13:00-18:45 (rows 54-77)
#    are set to 1, and rows for 20:00-21:45 (rows 82-89) revert to 0 (the block
#    that is "on" shifted 8 hours earlier in the day).
# 3) Update the sheet view (top-left cell / selection) to match the new state.
# 4) Apply (then effectively leave in place) an AutoFilter over A1:A97 so Excel
#    records the _xlnm._FilterDatabase defined name, matching the new window view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Shift dates in column A forward by one day ---
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 1
}

# --- 2) Update payload values in column B ---
# First reset the previously "on" block (rows 82-89) back to 0
for ($r = 82; $r -le 89; $r++) {
    $ws.Cells.Item($r, 2).Value2 = 0
}

# Then set the new "on" block (rows 54-77) to 1
for ($r = 54; $r -le 77; $r++) {
    $ws.Cells.Item($r, 2).Value2 = 1
}

# --- 3) Apply an AutoFilter over the Time Block column, producing the
#         _xlnm._FilterDatabase defined name seen in the workbook XML ---
$ws.Range("A1:A97").AutoFilter(1)

# --- 4) Update the visible window / selection to match the new state ---
$ws.Range("D66").Select()
$excel.ActiveWindow.ScrollRow = 62
